{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Target change: the paragraph containing the C++ snippet\n//   \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}};\"\n// gets extra point entries appended before the closing brace/semicolon, so\n// the final paragraph text becomes:\n//   \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}, {27,25}, {31,30}, {2,7}, {9,1));\"\n\nconst originalText =\n  \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}};\";\nconst newText =\n  \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}, {27,25}, {31,30}, {2,7}, {9,1));\";\n\nconst body = context.document.body;\nconst results = body.search(originalText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target paragraph text not found: \" + originalText);\n}\n\n// Replace the whole matched range (the entire original run's text) with the\n// new, longer text in one shot.\nresults.items[0].insertText(newText, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Target change: the paragraph containing the C++ snippet\n#   \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}};\"\n# gets extra point entries appended before the closing brace/semicolon, so\n# the final paragraph text becomes:\n#   \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}, {27,25}, {31,30}, {2,7}, {9,1));\"\n\n$d = $word.ActiveDocument\n\n$originalText = \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}};\"\n$newText = \"   poi P[] = {{4, 1}, {15, 20}, {30, 40}, {8, 4}, {13, 11}, {5, 6}, {27,25}, {31,30}, {2,7}, {9,1));\"\n\n$rng = $d.Content\n# Find.Execute(FindText, MatchCase, MatchWholeWord, MatchWildcards, ...)\n$found = $rng.Find.Execute($originalText, $true, $false, $false)\n\nif (-not $found) {\n    throw \"Target paragraph text not found: $originalText\"\n}\n\n$rng.Text = $newText\n"}
